$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("contacts")

# Update header row
$ws.Range("A1").Value = "Title"
$ws.Range("B1").Value = "Surname"
$ws.Range("C1").Value = "Lastname"

# Update "firstname" column values to their new contents
$ws.Range("B2").Value = "Lead"
$ws.Range("B3").Value = "Friend"
$ws.Range("B4").Value = "Lead"

# Reflect the updated selection saved in the workbook
$ws.Range("E6").Select() | Out-Null
